# IST price update 2025-12-21 10:46
# Insert a new column before column B (shifts existing B:M -> C:N) and
# populate it with the latest price-tracker scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at B; this pushes the existing B:M columns to C:N.
$ws.Range("B1").EntireColumn.Insert()

# Match the column width used by all the other price-history columns
# (column C, the former column B, still has the original width).
$ws.Range("B1").ColumnWidth = $ws.Range("C1").ColumnWidth

# New timestamp header for the freshly inserted column.
$ws.Range("B1").Value = "2025-12-21 16:12"

# Latest scraped prices for the new column (same price currently listed for
# each SKU at the time of this scrape).
$prices = @{
    2  = 929
    3  = 569
    4  = 299
    5  = 569
    6  = 499
    7  = 569
    8  = 929
    9  = 299
    10 = 299
    11 = 929
    12 = 569
    13 = 569
    14 = 499
    15 = 499
    16 = 299
    17 = 929
    18 = 499
    19 = 1497
    20 = 929
    21 = 499
    22 = 299
    23 = 1299
    24 = 929
    25 = 929
    26 = 1299
}

foreach ($row in $prices.Keys) {
    $ws.Cells.Item($row, 2).Value = $prices[$row]
}
